$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-06-05"

# Update the June row label in column A (row 7)
$ws.Range("A7").Value = "June (through 06-05)"

# Update June row (row 7) values for columns C, E, F, G, H, I
$ws.Range("C7").Value = 7
$ws.Range("E7").Value = 14
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 31
$ws.Range("H7").Value = 20
$ws.Range("I7").Value = 14

# Update Total row (row 8) values for columns C, E, F, G, H, I
$ws.Range("C8").Value = 216
$ws.Range("E8").Value = 309
$ws.Range("F8").Value = 210
$ws.Range("G8").Value = 389
$ws.Range("H8").Value = 651
$ws.Range("I8").Value = 678
